$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the missing "4920" values in column A for the rows that previously
# lacked them (these rows represent group/header rows for a new section).
# Column A is formatted as text (numFmtId 49 / "@"), so a direct numeric
# assignment would be stored as a text string. To store a genuine number
# while keeping the existing text-style (style index 2) applied to the
# cell, temporarily switch the cell to the default "Normal" style, assign
# the numeric value, then restore the original number format (copied from
# a neighbouring cell that already carries the desired style).
$rows = @(2, 11, 20, 29, 32, 33)
$formatSource = $ws.Cells.Item(3, 1)
foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Style = "Normal"
    $cell.Value = 4920
    $cell.NumberFormat = $formatSource.NumberFormat
}

# Update the active selection to A2:A33 with A2 as the active cell.
$ws.Range("A2:A33").Select()
